$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D:E").Insert()

$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 1074400
$ws.Cells.Item(8, 5).Value = 1037600
$ws.Cells.Item(9, 4).Value = 600100
$ws.Cells.Item(9, 5).Value = 577500
$ws.Cells.Item(10, 4).Value = 474300
$ws.Cells.Item(10, 5).Value = 460100
$ws.Cells.Item(12, 4).Value = 25300
$ws.Cells.Item(12, 5).Value = 23100
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 881000
$ws.Cells.Item(17, 5).Value = 833400
$ws.Cells.Item(18, 4).Value = 193400
$ws.Cells.Item(18, 5).Value = 204200
$ws.Cells.Item(20, 4).Value = 2400
$ws.Cells.Item(20, 5).Value = 2500
$ws.Cells.Item(21, 4).Value = 231500
$ws.Cells.Item(21, 5).Value = 242200
$ws.Cells.Item(22, 4).Value = 19700
$ws.Cells.Item(22, 5).Value = 19400
$ws.Cells.Item(23, 4).Value = 176100
$ws.Cells.Item(23, 5).Value = 187300
$ws.Cells.Item(24, 4).Value = 33300
$ws.Cells.Item(24, 5).Value = 41000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 142800
$ws.Cells.Item(26, 5).Value = 146300
$ws.Cells.Item(27, 4).Value = 142800
$ws.Cells.Item(27, 5).Value = 146300
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -2400
$ws.Cells.Item(32, 5).Value = -2500
$ws.Cells.Item(33, 4).Value = 142800
$ws.Cells.Item(33, 5).Value = 146300
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 142800
$ws.Cells.Item(35, 5).Value = 146300
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 316700
$ws.Cells.Item(41, 5).Value = 188300
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 345300
$ws.Cells.Item(43, 5).Value = 361500
$ws.Cells.Item(44, 4).Value = 382800
$ws.Cells.Item(44, 5).Value = 390600
$ws.Cells.Item(45, 4).Value = 33400
$ws.Cells.Item(45, 5).Value = 19900
$ws.Cells.Item(46, 4).Value = 1078200
$ws.Cells.Item(46, 5).Value = 960300
$ws.Cells.Item(47, 4).Value = 8500
$ws.Cells.Item(47, 5).Value = 8800
$ws.Cells.Item(48, 4).Value = 598200
$ws.Cells.Item(48, 5).Value = 587900
$ws.Cells.Item(49, 4).Value = 4266900
$ws.Cells.Item(49, 5).Value = 4287100
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 117400
$ws.Cells.Item(52, 5).Value = 123200
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 6069200
$ws.Cells.Item(54, 5).Value = 5967300
$ws.Cells.Item(57, 4).Value = 430200
$ws.Cells.Item(57, 5).Value = 449500
$ws.Cells.Item(58, 4).Value = 598300
$ws.Cells.Item(58, 5).Value = 301400
$ws.Cells.Item(59, 4).Value = 297800
$ws.Cells.Item(59, 5).Value = 276500
$ws.Cells.Item(60, 4).Value = 1326300
$ws.Cells.Item(60, 5).Value = 1027400
$ws.Cells.Item(61, 4).Value = 1508800
$ws.Cells.Item(61, 5).Value = 1803500
$ws.Cells.Item(62, 4).Value = 780300
$ws.Cells.Item(62, 5).Value = 783000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 3615400
$ws.Cells.Item(66, 5).Value = 3613900
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 3832600
$ws.Cells.Item(72, 5).Value = 3743500
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 2453800
$ws.Cells.Item(76, 5).Value = 2353400
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 142800
$ws.Cells.Item(81, 5).Value = 146300
$ws.Cells.Item(83, 4).Value = 35700
$ws.Cells.Item(83, 5).Value = 35500
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 195600
$ws.Cells.Item(89, 5).Value = 245300
$ws.Cells.Item(91, 4).Value = -30000
$ws.Cells.Item(91, 5).Value = -10800
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -29800
$ws.Cells.Item(94, 5).Value = -11100
$ws.Cells.Item(96, 4).Value = -53700
$ws.Cells.Item(96, 5).Value = -53400
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -34900
$ws.Cells.Item(100, 5).Value = -134600
$ws.Cells.Item(101, 4).Value = -2500
$ws.Cells.Item(101, 5).Value = -600
$ws.Cells.Item(102, 4).Value = 128400
$ws.Cells.Item(102, 5).Value = 99000
